# Auto-generated Excel COM-interop script to apply numeric corrections
# to the "Kujata_Profits" worksheets (recipe profit calculations).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1772.8636
$ws.Range("I40").Value = 1750.25
$ws.Range("K40").Value = 1750.25
$ws.Range("M40").Value = -1575.25
$ws.Range("H64").Value = 3822.111
$ws.Range("I64").Value = 4269.5
$ws.Range("J64").Value = 3694.2856
$ws.Range("K64").Value = 4269.5
$ws.Range("L64").Value = 3694.2856
$ws.Range("M64").Value = -4021.5
$ws.Range("N64").Value = -4190.2856
$ws.Range("H67").Value = 3822.111
$ws.Range("I67").Value = 4269.5
$ws.Range("J67").Value = 3694.2856
$ws.Range("K67").Value = 4269.5
$ws.Range("L67").Value = 3694.2856
$ws.Range("M67").Value = -3411.5
$ws.Range("N67").Value = -5410.2856
$ws.Range("H70").Value = 1749.5385
$ws.Range("I70").Value = 1740
$ws.Range("J70").Value = 1753.7778
$ws.Range("K70").Value = 5220
$ws.Range("L70").Value = 5261.3334
$ws.Range("M70").Value = -4950
$ws.Range("N70").Value = -5801.3334
$ws.Range("H73").Value = 1749.5385
$ws.Range("I73").Value = 1740
$ws.Range("J73").Value = 1753.7778
$ws.Range("K73").Value = 5220
$ws.Range("L73").Value = 5261.3334
$ws.Range("M73").Value = -4284
$ws.Range("N73").Value = -7133.3334
$ws.Range("H74").Value = 3241.9167
$ws.Range("I74").Value = 3343.2856
$ws.Range("J74").Value = 3100
$ws.Range("K74").Value = 3343.2856
$ws.Range("L74").Value = 3100
$ws.Range("M74").Value = -2407.2856
$ws.Range("N74").Value = -4972
$ws.Range("H77").Value = 3241.9167
$ws.Range("I77").Value = 3343.2856
$ws.Range("J77").Value = 3100
$ws.Range("K77").Value = 16716.428
$ws.Range("L77").Value = 15500
$ws.Range("M77").Value = -12036.428
$ws.Range("N77").Value = -24860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 8444
$ws.Range("J9").Value = 8444
$ws.Range("L9").Value = 8444
$ws.Range("N9").Value = -8784
$ws.Range("H20").Value = 8444
$ws.Range("J20").Value = 8444
$ws.Range("L20").Value = 8444
$ws.Range("N20").Value = -8984
$ws.Range("H23").Value = 75003.57000000001
$ws.Range("J23").Value = 71251.75
$ws.Range("L23").Value = 71251.75
$ws.Range("N23").Value = -71769.75
$ws.Range("H32").Value = 5710.76
$ws.Range("I32").Value = 5725.8125
$ws.Range("J32").Value = 5349.5
$ws.Range("K32").Value = 5725.8125
$ws.Range("L32").Value = 5349.5
$ws.Range("M32").Value = -5438.8125
$ws.Range("N32").Value = -5923.5
$ws.Range("H37").Value = 21766.666
$ws.Range("I37").Value = 18900
$ws.Range("J37").Value = 27500
$ws.Range("K37").Value = 18900
$ws.Range("L37").Value = 27500
$ws.Range("M37").Value = -18627
$ws.Range("N37").Value = -28046
$ws.Range("H55").Value = 37932.668
$ws.Range("J55").Value = 37932.668
$ws.Range("L55").Value = 37932.668
$ws.Range("N55").Value = -38562.668
$ws.Range("H63").Value = 2261.641
$ws.Range("I63").Value = 2114.963
$ws.Range("J63").Value = 2591.6667
$ws.Range("K63").Value = 2114.963
$ws.Range("L63").Value = 2591.6667
$ws.Range("M63").Value = -1428.963
$ws.Range("N63").Value = -3963.6667
$ws.Range("H66").Value = 2261.641
$ws.Range("I66").Value = 2114.963
$ws.Range("J66").Value = 2591.6667
$ws.Range("K66").Value = 10574.815
$ws.Range("L66").Value = 12958.3335
$ws.Range("M66").Value = -7142.815000000001
$ws.Range("N66").Value = -19822.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = ""
$ws.Range("N35").Value = 0
$ws.Range("H134").Value = 6302.8
$ws.Range("I134").Value = 873.8823
$ws.Range("J134").Value = 37066.668
$ws.Range("K134").Value = 2621.6469
$ws.Range("L134").Value = 111200.004
$ws.Range("M134").Value = -86.64689999999973
$ws.Range("N134").Value = -116270.004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6062951.5
$ws.Range("I62").Value = 2396.6667
$ws.Range("J62").Value = 66668500
$ws.Range("K62").Value = 2396.6667
$ws.Range("L62").Value = 66668500
$ws.Range("M62").Value = -1772.6667
$ws.Range("N62").Value = -66669748
$ws.Range("H65").Value = 6062951.5
$ws.Range("I65").Value = 2396.6667
$ws.Range("J65").Value = 66668500
$ws.Range("K65").Value = 11983.3335
$ws.Range("L65").Value = 333342500
$ws.Range("M65").Value = -8863.333500000001
$ws.Range("N65").Value = -333348740
$ws.Range("H109").Value = 12967
$ws.Range("J109").Value = 12967
$ws.Range("L109").Value = 12967
$ws.Range("N109").Value = -15047
$ws.Range("H122").Value = 4252.9644
$ws.Range("I122").Value = 4488.1924
$ws.Range("K122").Value = 13464.5772
$ws.Range("M122").Value = -11014.5772
$ws.Range("H132").Value = 2491.2942
$ws.Range("I132").Value = 2156.9333
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6470.7999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3940.7999
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 14287466
$ws.Range("I134").Value = 1859.1724
$ws.Range("J134").Value = 83334570
$ws.Range("K134").Value = 5577.5172
$ws.Range("L134").Value = 250003710
$ws.Range("M134").Value = -3042.5172
$ws.Range("N134").Value = -250008780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 9301.182000000001
$ws.Range("J82").Value = 12000
$ws.Range("L82").Value = 36000
$ws.Range("N82").Value = -36812
$ws.Range("H85").Value = 9301.182000000001
$ws.Range("J85").Value = 12000
$ws.Range("L85").Value = 36000
$ws.Range("N85").Value = -38808
$ws.Range("H113").Value = 683.03845
$ws.Range("J113").Value = 717.6842
$ws.Range("L113").Value = 2153.0526
$ws.Range("N113").Value = -6493.0526
$ws.Range("H131").Value = 26317438
$ws.Range("J131").Value = 1806.1471
$ws.Range("L131").Value = 5418.4413
$ws.Range("N131").Value = -15498.4413

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 15000.111
$ws.Range("J63").Value = 15000.111
$ws.Range("L63").Value = 15000.111
$ws.Range("N63").Value = -16372.111
$ws.Range("H66").Value = 15000.111
$ws.Range("J66").Value = 15000.111
$ws.Range("L66").Value = 45000.333
$ws.Range("N66").Value = -51864.333
$ws.Range("H102").Value = 795.63635
$ws.Range("I102").Value = 609.8
$ws.Range("J102").Value = 1193.8572
$ws.Range("K102").Value = 609.8
$ws.Range("L102").Value = 1193.8572
$ws.Range("M102").Value = 1012.2
$ws.Range("N102").Value = -4437.8572
$ws.Range("H132").Value = 4857.0835
$ws.Range("I132").Value = 6070.5
$ws.Range("K132").Value = 18211.5
$ws.Range("M132").Value = -15681.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2990.0908
$ws.Range("I40").Value = 2004.3158
$ws.Range("K40").Value = 2004.3158
$ws.Range("M40").Value = -1868.3158
$ws.Range("H56").Value = 11333.333
$ws.Range("H68").Value = 2990
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 2990
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H112").Value = 48499
$ws.Range("J112").Value = 48499
$ws.Range("L112").Value = 48499
$ws.Range("N112").Value = -51453
$ws.Range("H122").Value = 31251924
$ws.Range("H132").Value = 2797.55
$ws.Range("I132").Value = 2410.9285
$ws.Range("J132").Value = 3699.6667
$ws.Range("K132").Value = 7232.7855
$ws.Range("L132").Value = 11099.0001
$ws.Range("M132").Value = -4702.7855
$ws.Range("N132").Value = -16159.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -710
$ws.Range("H61").Value = 11057
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("H62").Value = 166674160
$ws.Range("J62").Value = 11251.5
$ws.Range("L62").Value = 11251.5
$ws.Range("N62").Value = -12499.5
$ws.Range("H65").Value = 166674160
$ws.Range("J65").Value = 11251.5
$ws.Range("L65").Value = 56257.5
$ws.Range("N65").Value = -62497.5
$ws.Range("H109").Value = 36039.8
$ws.Range("J109").Value = 32464.25
$ws.Range("L109").Value = 32464.25
$ws.Range("N109").Value = -35238.25
$ws.Range("H122").Value = 19232534
$ws.Range("I122").Value = 25001992
$ws.Range("K122").Value = 75005976
$ws.Range("M122").Value = -75003526
$ws.Range("H132").Value = 1358.0435
$ws.Range("I132").Value = 921.85
$ws.Range("J132").Value = 4266
$ws.Range("K132").Value = 2765.55
$ws.Range("L132").Value = 12798
$ws.Range("M132").Value = -235.5500000000002
$ws.Range("N132").Value = -17858

